$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily auction rows appended to the bottom of the table (MV update).
$data = @(
    @("21-09-2021", 40, 202, 40, 785),
    @("22-09-2021", 40, 150, 40, 785),
    @("23-09-2021", 40, 104, 40, 788),
    @("24-09-2021", 40, 152, 40, 792),
    @("27-09-2021", 40, 166, 40, 796),
    @("28-09-2021", 40, 125, 40, 799),
    @("29-09-2021", 40, 155, 40, 808),
    @("30-09-2021", 40, 151, 40, 811),
    @("01-10-2021", 40, 143, 40, 803)
)

$startRow = 172
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $dateText = $row[0]
    $dateCell = $ws.Cells.Item($r, 1)

    # "01-10-2021" is a valid dd-mm-yyyy *and* mm-dd-yyyy date, so Excel's
    # automatic type detection would silently turn it into a date serial
    # number. Format the cell as text first so the literal string (matching
    # the other "dd-mm-yyyy" labels already in column A) is preserved, then
    # drop back to the default style so no visible formatting changes.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
